# Update NATMI LR-pair stats (Inhba-Bambi) per Dr Hou advice: ligand/receptor-expressing
# cell counts go from 1 to 3 per group, with recomputed average/total expression,
# specificity, and edge-weight statistics for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.288150666666667
$ws.Range("H2").Value = 6.864452
$ws.Range("I2").Value = 0.3964219041944151
$ws.Range("J2").Value = 0.3964219041944151
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.201683333333333
$ws.Range("N2").Value = 9.60505
$ws.Range("O2").Value = 0.5163673346595562
$ws.Range("P2").Value = 0.5163673346595563
$ws.Range("Q2").Value = 7.325933853622223
$ws.Range("R2").Value = 65.9334046826
$ws.Range("S2").Value = 0.2046993220695361
$ws.Range("T2").Value = 0.2046993220695361

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.288150666666667
$ws.Range("H3").Value = 6.864452
$ws.Range("I3").Value = 0.3964219041944151
$ws.Range("J3").Value = 0.3964219041944151
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9983063333333334
$ws.Range("N3").Value = 2.994919
$ws.Range("O3").Value = 0.1610067976274214
$ws.Range("P3").Value = 0.1610067976274214
$ws.Range("Q3").Value = 2.284275302154223
$ws.Range("R3").Value = 20.558477719388
$ws.Range("S3").Value = 0.06382662130370721
$ws.Range("T3").Value = 0.06382662130370721

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.288150666666667
$ws.Range("H4").Value = 6.864452
$ws.Range("I4").Value = 0.3964219041944151
$ws.Range("J4").Value = 0.3964219041944151
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.562824666666667
$ws.Range("N4").Value = 4.688474
$ws.Range("O4").Value = 0.2520522873905527
$ws.Range("P4").Value = 0.2520522873905527
$ws.Range("Q4").Value = 3.575978302916444
$ws.Range("R4").Value = 32.183804726248
$ws.Range("S4").Value = 0.09991904772392088
$ws.Range("T4").Value = 0.09991904772392088

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.288150666666667
$ws.Range("H5").Value = 6.864452
$ws.Range("I5").Value = 0.3964219041944151
$ws.Range("J5").Value = 0.3964219041944151
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4375843333333334
$ws.Range("N5").Value = 1.312753
$ws.Range("O5").Value = 0.07057358032246958
$ws.Range("P5").Value = 0.07057358032246959
$ws.Range("Q5").Value = 1.001258884039556
$ws.Range("R5").Value = 9.011329956356
$ws.Range("S5").Value = 0.02797691309725089
$ws.Range("T5").Value = 0.0279769130972509

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.588894
$ws.Range("H6").Value = 7.766681999999999
$ws.Range("I6").Value = 0.4485256605643812
$ws.Range("J6").Value = 0.4485256605643813
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.201683333333333
$ws.Range("N6").Value = 9.60505
$ws.Range("O6").Value = 0.5163673346595562
$ws.Range("P6").Value = 0.5163673346595563
$ws.Range("Q6").Value = 8.288818771566666
$ws.Range("R6").Value = 74.5993689441
$ws.Range("S6").Value = 0.2316039998720463
$ws.Range("T6").Value = 0.2316039998720464

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.588894
$ws.Range("H7").Value = 7.766681999999999
$ws.Range("I7").Value = 0.4485256605643812
$ws.Range("J7").Value = 0.4485256605643813
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9983063333333334
$ws.Range("N7").Value = 2.994919
$ws.Range("O7").Value = 0.1610067976274214
$ws.Range("P7").Value = 0.1610067976274214
$ws.Range("Q7").Value = 2.584509276528667
$ws.Range("R7").Value = 23.260583488758
$ws.Range("S7").Value = 0.07221568026119482
$ws.Range("T7").Value = 0.07221568026119482

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.588894
$ws.Range("H8").Value = 7.766681999999999
$ws.Range("I8").Value = 0.4485256605643812
$ws.Range("J8").Value = 0.4485256605643813
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.562824666666667
$ws.Range("N8").Value = 4.688474
$ws.Range("O8").Value = 0.2520522873905527
$ws.Range("P8").Value = 0.2520522873905527
$ws.Range("Q8").Value = 4.045987402585333
$ws.Range("R8").Value = 36.413886623268
$ws.Range("S8").Value = 0.1130519186986109
$ws.Range("T8").Value = 0.1130519186986109

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.588894
$ws.Range("H9").Value = 7.766681999999999
$ws.Range("I9").Value = 0.4485256605643812
$ws.Range("J9").Value = 0.4485256605643813
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4375843333333334
$ws.Range("N9").Value = 1.312753
$ws.Range("O9").Value = 0.07057358032246958
$ws.Range("P9").Value = 0.07057358032246959
$ws.Range("Q9").Value = 1.132859455060667
$ws.Range("R9").Value = 10.195735095546
$ws.Range("S9").Value = 0.03165406173252908
$ws.Range("T9").Value = 0.03165406173252909

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.3123523333333333
$ws.Range("H10").Value = 0.9370569999999999
$ws.Range("I10").Value = 0.05411501461132016
$ws.Range("J10").Value = 0.05411501461132018
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 3.201683333333333
$ws.Range("N10").Value = 9.60505
$ws.Range("O10").Value = 0.5163673346595562
$ws.Range("P10").Value = 0.5163673346595563
$ws.Range("Q10").Value = 1.000053259761111
$ws.Range("R10").Value = 9.000479337849999
$ws.Range("S10").Value = 0.02794322585991033
$ws.Range("T10").Value = 0.02794322585991035

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.3123523333333333
$ws.Range("H11").Value = 0.9370569999999999
$ws.Range("I11").Value = 0.05411501461132016
$ws.Range("J11").Value = 0.05411501461132018
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.9983063333333334
$ws.Range("N11").Value = 2.994919
$ws.Range("O11").Value = 0.1610067976274214
$ws.Range("P11").Value = 0.1610067976274214
$ws.Range("Q11").Value = 0.3118233125981111
$ws.Range("R11").Value = 2.806409813383
$ws.Range("S11").Value = 0.008712885206129777
$ws.Range("T11").Value = 0.008712885206129778

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.3123523333333333
$ws.Range("H12").Value = 0.9370569999999999
$ws.Range("I12").Value = 0.05411501461132016
$ws.Range("J12").Value = 0.05411501461132018
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1.562824666666667
$ws.Range("N12").Value = 4.688474
$ws.Range("O12").Value = 0.2520522873905527
$ws.Range("P12").Value = 0.2520522873905527
$ws.Range("Q12").Value = 0.4881519312242222
$ws.Range("R12").Value = 4.393367381018
$ws.Range("S12").Value = 0.01363981321495643
$ws.Range("T12").Value = 0.01363981321495643

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.3123523333333333
$ws.Range("H13").Value = 0.9370569999999999
$ws.Range("I13").Value = 0.05411501461132016
$ws.Range("J13").Value = 0.05411501461132018
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4375843333333334
$ws.Range("N13").Value = 1.312753
$ws.Range("O13").Value = 0.07057358032246958
$ws.Range("P13").Value = 0.07057358032246959
$ws.Range("Q13").Value = 0.1366804875467778
$ws.Range("R13").Value = 1.230124387921
$ws.Range("S13").Value = 0.003819090330323618
$ws.Range("T13").Value = 0.00381909033032362

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.5826116666666666
$ws.Range("H14").Value = 1.747835
$ws.Range("I14").Value = 0.1009374206298835
$ws.Range("J14").Value = 0.1009374206298836
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3.201683333333333
$ws.Range("N14").Value = 9.60505
$ws.Range("O14").Value = 0.5163673346595562
$ws.Range("P14").Value = 0.5163673346595563
$ws.Range("Q14").Value = 1.865338062972222
$ws.Range("R14").Value = 16.78804256675
$ws.Range("S14").Value = 0.05212078685806347
$ws.Range("T14").Value = 0.05212078685806349

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.5826116666666666
$ws.Range("H15").Value = 1.747835
$ws.Range("I15").Value = 0.1009374206298835
$ws.Range("J15").Value = 0.1009374206298836
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.9983063333333334
$ws.Range("N15").Value = 2.994919
$ws.Range("O15").Value = 0.1610067976274214
$ws.Range("P15").Value = 0.1610067976274214
$ws.Range("Q15").Value = 0.5816249167072223
$ws.Range("R15").Value = 5.234624250365
$ws.Range("S15").Value = 0.01625161085638957
$ws.Range("T15").Value = 0.01625161085638957

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.5826116666666666
$ws.Range("H16").Value = 1.747835
$ws.Range("I16").Value = 0.1009374206298835
$ws.Range("J16").Value = 0.1009374206298836
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 1.562824666666667
$ws.Range("N16").Value = 4.688474
$ws.Range("O16").Value = 0.2520522873905527
$ws.Range("P16").Value = 0.2520522873905527
$ws.Range("Q16").Value = 0.9105198837544445
$ws.Range("R16").Value = 8.19467895379
$ws.Range("S16").Value = 0.02544150775306452
$ws.Range("T16").Value = 0.02544150775306452

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.5826116666666666
$ws.Range("H17").Value = 1.747835
$ws.Range("I17").Value = 0.1009374206298835
$ws.Range("J17").Value = 0.1009374206298836
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.4375843333333334
$ws.Range("N17").Value = 1.312753
$ws.Range("O17").Value = 0.07057358032246958
$ws.Range("P17").Value = 0.07057358032246959
$ws.Range("Q17").Value = 0.2549417377505556
$ws.Range("R17").Value = 2.294475639755
$ws.Range("S17").Value = 0.007123515162365984
$ws.Range("T17").Value = 0.007123515162365987

